$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: confusion_matrix (inline strings, multi-line text)
$ws.Range("B2").Value = "[[96503  4301]`n [ 2516    53]]"
$ws.Range("C2").Value = "[[98920  1884]`n [ 2529    40]]"
$ws.Range("D2").Value = "[[100403    401]`n [   560   2009]]"
$ws.Range("E2").Value = "[[95738  5066]`n [  560  2009]]"

# Row 3: accuracy_score (numeric)
$ws.Range("B3").Value = 0.9340543468797462
$ws.Range("C3").Value = 0.9573099358633299
$ws.Range("D3").Value = 0.9907035686301066
$ws.Range("E3").Value = 0.9455757306066381

# Row 4: f1_score (numeric)
$ws.Range("B4").Value = 0.01531128123645818
$ws.Range("C4").Value = 0.01780547518361896
$ws.Range("D4").Value = 0.8069893552922273
$ws.Range("E4").Value = 0.4166321028618831
